# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures scraped by the scheduled runner
# to the Famfrit_Profits workbook (one worksheet per crafting job).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 167
$ws.Range("I8").Value = 151.91667
$ws.Range("K8").Value = 455.75001
$ws.Range("M8").Value = -316.75001
$ws.Range("H92").Value = 952.84
$ws.Range("I92").Value = 774.86957
$ws.Range("J92").Value = 2999.5
$ws.Range("K92").Value = 774.86957
$ws.Range("L92").Value = 2999.5
$ws.Range("M92").Value = 473.13043
$ws.Range("N92").Value = -5495.5
$ws.Range("H95").Value = 43437.668
$ws.Range("J95").Value = 43437.668
$ws.Range("L95").Value = 43437.668
$ws.Range("N95").Value = -48929.668
$ws.Range("H100").Value = 2200.45
$ws.Range("I100").Value = 1833.1818
$ws.Range("J100").Value = 2649.3333
$ws.Range("K100").Value = 1833.1818
$ws.Range("L100").Value = 2649.3333
$ws.Range("M100").Value = -1292.1818
$ws.Range("N100").Value = -3731.3333
$ws.Range("H112").Value = 3029.9614
$ws.Range("I112").Value = 948.5454999999999
$ws.Range("K112").Value = 2845.6365
$ws.Range("M112").Value = -1737.6365
$ws.Range("H116").Value = 4050.353
$ws.Range("I116").Value = 3673.7
$ws.Range("J116").Value = 4588.4287
$ws.Range("K116").Value = 3673.7
$ws.Range("L116").Value = 4588.4287
$ws.Range("M116").Value = -231.6999999999998
$ws.Range("N116").Value = -11472.4287
$ws.Range("H132").Value = 8006.278
$ws.Range("I132").Value = 7947.8823
$ws.Range("K132").Value = 23843.6469
$ws.Range("M132").Value = -21313.6469
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 3451.761
$ws.Range("I137").Value = 2335.6177
$ws.Range("J137").Value = 6614.1665
$ws.Range("K137").Value = 7006.853099999999
$ws.Range("L137").Value = 19842.4995
$ws.Range("M137").Value = -4456.853099999999
$ws.Range("N137").Value = -24942.4995
$ws.Range("H138").Value = 13000.611
$ws.Range("I138").Value = 4565.3335
$ws.Range("J138").Value = 14687.667
$ws.Range("K138").Value = 13696.0005
$ws.Range("L138").Value = 44063.001
$ws.Range("M138").Value = -8556.000499999998
$ws.Range("N138").Value = -54343.001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19233960
$ws.Range("I32").Value = 21278374
$ws.Range("J32").Value = 16478.6
$ws.Range("K32").Value = 21278374
$ws.Range("L32").Value = 16478.6
$ws.Range("M32").Value = -21278087
$ws.Range("N32").Value = -17052.6
$ws.Range("H61").Value = 71431170
$ws.Range("I61").Value = 125001300
$ws.Range("K61").Value = 125001300
$ws.Range("M61").Value = -125001088
$ws.Range("H97").Value = 1491.8636
$ws.Range("I97").Value = 1491.8636
$ws.Range("K97").Value = 1491.8636
$ws.Range("M97").Value = -995.8635999999999
$ws.Range("H122").Value = 14495309
$ws.Range("I122").Value = 2322.0952
$ws.Range("K122").Value = 6966.285600000001
$ws.Range("M122").Value = -4516.285600000001
$ws.Range("H132").Value = 43553440
$ws.Range("I132").Value = 20409
$ws.Range("J132").Value = 125177880
$ws.Range("K132").Value = 61227
$ws.Range("L132").Value = 375533640
$ws.Range("M132").Value = -58697
$ws.Range("N132").Value = -375538700
$ws.Range("H136").Value = 71431170
$ws.Range("I136").Value = 125001300
$ws.Range("K136").Value = 375003900
$ws.Range("M136").Value = -375001350

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 355.3
$ws.Range("J80").Value = 72
$ws.Range("L80").Value = 72
$ws.Range("N80").Value = -2068
$ws.Range("H83").Value = 355.3
$ws.Range("J83").Value = 72
$ws.Range("L83").Value = 360
$ws.Range("N83").Value = -10344
$ws.Range("H99").Value = 4619
$ws.Range("I99").Value = 2766.3333
$ws.Range("K99").Value = 2766.3333
$ws.Range("M99").Value = -1268.3333

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1001
$ws.Range("I3").Value = 999
$ws.Range("J3").Value = 1003
$ws.Range("K3").Value = 999
$ws.Range("L3").Value = 1003
$ws.Range("M3").Value = -886
$ws.Range("N3").Value = -1229
$ws.Range("H120").Value = 26902
$ws.Range("J120").Value = 24663.8
$ws.Range("L120").Value = 24663.8
$ws.Range("N120").Value = -31921.8
$ws.Range("H132").Value = 86195.414
$ws.Range("H134").Value = 5397.1665
$ws.Range("I134").Value = 5397.1665
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16191.4995
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13656.4995
$ws.Range("N134").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 32155416
$ws.Range("I4").Value = 40493444
$ws.Range("J4").Value = 12700017
$ws.Range("K4").Value = 121480332
$ws.Range("L4").Value = 38100051
$ws.Range("M4").Value = -121480220
$ws.Range("N4").Value = -38100275
$ws.Range("H6").Value = 78.666664
$ws.Range("I6").Value = 78.666664
$ws.Range("K6").Value = 235.999992
$ws.Range("M6").Value = -122.999992
$ws.Range("H70").Value = 2698.8
$ws.Range("I70").Value = 1831.3334
$ws.Range("K70").Value = 5494.0002
$ws.Range("M70").Value = -5179.0002
$ws.Range("H73").Value = 2698.8
$ws.Range("I73").Value = 1831.3334
$ws.Range("K73").Value = 5494.0002
$ws.Range("M73").Value = -4402.0002
$ws.Range("H107").Value = 1887.3
$ws.Range("I107").Value = 548.5
$ws.Range("K107").Value = 1645.5
$ws.Range("M107").Value = 274.5
$ws.Range("H114").Value = 2387.3333
$ws.Range("J114").Value = 3531
$ws.Range("L114").Value = 10593
$ws.Range("N114").Value = -17101
$ws.Range("H134").Value = 5185
$ws.Range("I134").Value = 1491.6666
$ws.Range("K134").Value = 4474.9998
$ws.Range("M134").Value = 595.0002000000004
$ws.Range("H138").Value = 1170.6
$ws.Range("I138").Value = 1170.6
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 3511.8
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1628.2
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 2800.28
$ws.Range("I139").Value = 1543.7273
$ws.Range("J139").Value = 12015
$ws.Range("K139").Value = 4631.1819
$ws.Range("L139").Value = 36045
$ws.Range("M139").Value = 508.8181000000004
$ws.Range("N139").Value = -46325

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4196.8486
$ws.Range("I132").Value = 3733.8262
$ws.Range("J132").Value = 5261.8
$ws.Range("K132").Value = 11201.4786
$ws.Range("L132").Value = 15785.4
$ws.Range("M132").Value = -8671.4786
$ws.Range("N132").Value = -20845.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 372197.03
$ws.Range("I93").Value = 2462.6206
$ws.Range("J93").Value = 1265721.9
$ws.Range("K93").Value = 2462.6206
$ws.Range("L93").Value = 1265721.9
$ws.Range("M93").Value = -1214.6206
$ws.Range("N93").Value = -1268217.9
$ws.Range("H122").Value = 17861528
$ws.Range("I122").Value = 3134.6667
$ws.Range("J122").Value = 31255324
$ws.Range("K122").Value = 9404.000100000001
$ws.Range("L122").Value = 93765972
$ws.Range("M122").Value = -6954.000100000001
$ws.Range("N122").Value = -93770872
$ws.Range("H132").Value = 133336800
$ws.Range("I132").Value = 3540
$ws.Range("J132").Value = 500003260
$ws.Range("K132").Value = 10620
$ws.Range("L132").Value = 1500009780
$ws.Range("M132").Value = -8090
$ws.Range("N132").Value = -1500014840

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H14").Value = 2710.25
$ws.Range("I14").Value = 1946.091
$ws.Range("J14").Value = 4391.4
$ws.Range("K14").Value = 1946.091
$ws.Range("L14").Value = 4391.4
$ws.Range("M14").Value = -1778.091
$ws.Range("N14").Value = -4727.4
$ws.Range("H70").Value = 28666
$ws.Range("J70").Value = 28666
$ws.Range("L70").Value = 28666
$ws.Range("N70").Value = -29296
$ws.Range("H73").Value = 28666
$ws.Range("J73").Value = 28666
$ws.Range("L73").Value = 28666
$ws.Range("N73").Value = -30850
$ws.Range("H126").Value = 4675.5625
$ws.Range("I126").Value = 4131.759
$ws.Range("K126").Value = 12395.277
$ws.Range("M126").Value = -9925.277
$ws.Range("H132").Value = 1629.25
$ws.Range("I132").Value = 1589.1578
$ws.Range("J132").Value = 1781.6
$ws.Range("K132").Value = 4767.4734
$ws.Range("L132").Value = 5344.799999999999
$ws.Range("M132").Value = -2237.4734
$ws.Range("N132").Value = -10404.8
